$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Paragraph: "Desde o início dos anos 90 ... (décadas de 50, 60, 70, 80, 90,
# 2000, 2010 e 2016)." -- add the mustard/accent4 color to the paragraph mark
# and every run, and remove the stray _GoBack bookmark that was splitting the
# word "habilidades exigidas" into two runs ("habilidades e" / "xigidas...").
# ---------------------------------------------------------------------------
$targetText = "Desde o in"
$found = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith($targetText)) {
        $xmlPara = '<w:p ' + $wNs + ' w:rsidR="00937926" w:rsidRDefault="005F286A">' + `
            '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:hanging="360"/><w:contextualSpacing/>' + `
            '<w:rPr><w:color w:val="806000" w:themeColor="accent4" w:themeShade="80"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' + `
            '<w:r><w:rPr><w:color w:val="806000" w:themeColor="accent4" w:themeShade="80"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' + `
            '<w:t xml:space="preserve">Desde o início dos anos 90 a internet vem modificando o comportamento humano e, consequentemente, a forma como lidamos as tecnologias. Desenvolver </w:t></w:r>' + `
            '<w:proofErr w:type="spellStart"/>' + `
            '<w:r><w:rPr><w:color w:val="806000" w:themeColor="accent4" w:themeShade="80"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>software</w:t></w:r>' + `
            '<w:proofErr w:type="spellEnd"/>' + `
            '<w:r><w:rPr><w:color w:val="806000" w:themeColor="accent4" w:themeShade="80"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' + `
            '<w:t xml:space="preserve"> tornou-se então mais que um desafio de criar soluções para produtividade e automação, sendo que hoje as pessoas utilizam software para comunicação, pesquisar por produtos mais recomendados, fazer dispositivos voarem, dar aula, etc. Pensando assim, construa uma linha do tempo para definir qual o perfil de habilidades exigidas para os profissionais desde a década 50 até hoje (décadas de 50, 60, 70, 80, 90, 2000, 2010 e 2016).</w:t></w:r>' + `
            '</w:p>'
        $p.Range.InsertXML($xmlPara)
        $found = $true
        break
    }
}
if (-not $found) {
    Write-Host "WARNING: paragraph starting with 'Desde o in' was not found"
}

# ---------------------------------------------------------------------------
# Paragraph: "O que é um Processo de Software?" -- append a _GoBack bookmark
# at the end of the paragraph (after the existing run).
# ---------------------------------------------------------------------------
$targetText2 = "O que"
$found2 = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd().StartsWith($targetText2)) {
        $xmlPara = '<w:p ' + $wNs + ' w:rsidR="00937926" w:rsidRPr="00C97E24" w:rsidRDefault="005F286A">' + `
            '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:hanging="360"/><w:contextualSpacing/>' + `
            '<w:rPr><w:color w:val="00B050"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' + `
            '<w:r w:rsidRPr="00C97E24"><w:rPr><w:color w:val="00B050"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>' + `
            '<w:t>O que é um Processo de Software?</w:t></w:r>' + `
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
            '</w:p>'
        $p.Range.InsertXML($xmlPara)
        $found2 = $true
        break
    }
}
if (-not $found2) {
    Write-Host "WARNING: paragraph 'O que é um Processo de Software?' was not found"
}

Write-Host "edit complete"
